# Apply the MONTANA_2022 data-cleaning fixes:
#  1. Rename header columns to short machine-friendly names.
#  2. Normalize capitalization of connector words ("de"/"del"/"el"/"los")
#     to title case in a handful of place names.
#  3. Remove the trailing metadata/footer rows (133-137), shrinking the
#     used range down to A1:D131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -----------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Capitalization fixes in place names -----------------------------
$ws.Range("A27").Value = "Ciudad De México"
$ws.Range("A35").Value = "Estado De México"
$ws.Range("B36").Value = "San Felipe Del Progreso"
$ws.Range("B42").Value = "Apaseo El Alto"
$ws.Range("B49").Value = "San Luis De La Paz"
$ws.Range("B51").Value = "Chilpancingo De Los Bravo"
$ws.Range("B55").Value = "Huitzuco De Los Figueroa"
$ws.Range("B62").Value = "Pachuca De Soto"
$ws.Range("B65").Value = "Autlán De Navarro"
$ws.Range("B71").Value = "Zapotlán El Grande"
$ws.Range("B80").Value = "San Nicolás De Los Garza"
$ws.Range("B82").Value = "Oaxaca De Juárez"
$ws.Range("B83").Value = "Pinotepa De Don Luis"
$ws.Range("B84").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B88").Value = "Tepexi De Rodríguez"
$ws.Range("B90").Value = "Amealco De Bonfil"
$ws.Range("B95").Value = "Mexquitic De Carmona"
$ws.Range("B98").Value = "Villa De Ramos"
$ws.Range("B126").Value = "Concepción Del Oro"

# --- 3. Remove trailing metadata rows (133-137) -------------------------
$ws.Range("A133:D137").EntireRow.Delete()
